$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8077558279037476
$ws.Range("B1").Value = 1.259222507476807
$ws.Range("C1").Value = 4.532760143280029
$ws.Range("D1").Value = 4.09904146194458
$ws.Range("E1").Value = 0.8175625205039978
